$wb = $excel.ActiveWorkbook

# Update status text on the Overview sheet (row for 36e7dab6-... file):
# "Ready for handoff" -> "Handback transform failed"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B7").Value = "Handback transform failed"
$wsOverview.Range("C7").Value = "Handback transform failed"

# Add Error Detail (column K) for zh-cn sheet, row 7
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K7").Value = "Handback file name: 0pjs0wy3.xjd is different with handoff file name: 36e7dab6-6ad8-4f94-ae6f-6b12fed0cda4.473fb19fdd33313d689175ef688cb0ff5e5d9bd2.zh-cn."

# Add Error Detail (column K) for de-de sheet, row 7
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K7").Value = "Handback file name: 0pjs0wy3.xjd is different with handoff file name: 36e7dab6-6ad8-4f94-ae6f-6b12fed0cda4.473fb19fdd33313d689175ef688cb0ff5e5d9bd2.de-de."
